# Applies the commit "Enlace para ver diagrama er":
#  1. Adds "Diagrama ER: <link>" (as two separate runs) to the previously
#     empty "List Paragraph" that follows "Diagrama entidad relación." and
#     inserts a new empty "List Paragraph" right after it.
#  2. Removes the stray <w:lastRenderedPageBreak/> from the run that holds
#     "Plan de trabajo propuesto con tiempos y recursos."

$d = $word.ActiveDocument

# --- Helper XML package wrapper -------------------------------------------------
function New-WordPkg([string]$bodyInner) {
    return @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>$bodyInner<w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
}

# --- 1. Find the empty "Prrafodelista" paragraph right after ------------------
#        "Diagrama entidad relación." and fill it with the link text.
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "Diagrama entidad relación.`r") {
        $target = $d.Paragraphs.Item($i + 1)
        $found = $true
        break
    }
}

if (-not $found) {
    throw "Could not locate the 'Diagrama entidad relación.' paragraph"
}

$targetIndex = $i + 1
$insertPos = $target.Range.Start
$insertHere = $d.Range($insertPos, $insertPos)
$runsXml = '<w:p><w:r><w:t xml:space="preserve">Diagrama ER: </w:t></w:r><w:r><w:t>https://github.com/AndrewIbanhez/PruebaGH/blob/main/diagramaer.png</w:t></w:r></w:p>'
$insertHere.InsertXML((New-WordPkg $runsXml))

# --- 2. Insert a brand-new empty "Prrafodelista" paragraph right after --------
$filled = $d.Paragraphs.Item($targetIndex)
$afterPos = $filled.Range.End
$afterHere = $d.Range($afterPos, $afterPos)
$emptyParaXml = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/></w:pPr></w:p>'
$afterHere.InsertXML((New-WordPkg $emptyParaXml))

# --- 3. Remove <w:lastRenderedPageBreak/> from the "Plan de trabajo          --
#        propuesto con tiempos y recursos." run, keeping paragraph props.
$found2 = $false
for ($j = 1; $j -le $d.Paragraphs.Count; $j++) {
    $cand2 = $d.Paragraphs.Item($j)
    if ($cand2.Range.Text -eq "Plan de trabajo propuesto con tiempos y recursos.`r") {
        $planPara = $cand2
        $found2 = $true
        break
    }
}

if (-not $found2) {
    throw "Could not locate the 'Plan de trabajo propuesto...' paragraph"
}

$pRange = $planPara.Range
$runStart = $pRange.Start
$runEnd = $pRange.End - 1
$runRange = $d.Range($runStart, $runEnd)
$cleanRunXml = '<w:p><w:r><w:t>Plan de trabajo propuesto con tiempos y recursos.</w:t></w:r></w:p>'
$runRange.InsertXML((New-WordPkg $cleanRunXml))

Write-Output "Edit complete"
